$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "25.862.10"
$ws.Range("E2").Value = "  +0.24%  "

$ws.Range("D3").Value = "1.637.33"
$ws.Range("E3").Value = "  +0.74%  "

$ws.Range("D5").Value = "'215.42"
$ws.Range("E5").Value = "  +0.19%  "

$ws.Range("E6").Value = "  -0.47%  "

$ws.Range("E7").Value = "  +0.04%  "

$ws.Range("D8").Value = "'0.260"
$ws.Range("E8").Value = "  +0.63%  "

$ws.Range("D9").Value = "'0.0645"
$ws.Range("E9").Value = "  +1.15%  "

$ws.Range("D10").Value = "'20.23"
$ws.Range("E10").Value = "  +4.58%  "

$ws.Range("D11").Value = "'0.0781"
$ws.Range("E11").Value = "  +0.52%  "

$ws.Range("D12").Value = "1.657.81"
$ws.Range("E12").Value = "  +2.02%  "

$ws.Range("E13").Value = "  +0.34%  "

$ws.Range("D14").Value = "1.863.07"
$ws.Range("E14").Value = "  +0.78%  "

$ws.Range("D16").Value = "0.0₃0768"
$ws.Range("E16").Value = "  +2.07%  "

$ws.Range("D17").Value = "'63.32"
$ws.Range("E17").Value = "  -0.41%  "

$ws.Range("D18").Value = "25.873.57"
$ws.Range("E18").Value = "  +0.39%  "

$ws.Range("E19").Value = "  -0.04%  "

$ws.Range("D20").Value = "'194.46"
$ws.Range("E20").Value = "  +0.39%  "

$ws.Range("D21").Value = "'4.38"
$ws.Range("E21").Value = "  +1.28%  "

$ws.Range("D22").Value = "'9.97"
$ws.Range("E22").Value = "  +1.86%  "

$ws.Range("E23").Value = "  +3.88%  "

$ws.Range("E24").Value = "  +0.06%  "

$ws.Range("D25").Value = "'1.75"
$ws.Range("E25").Value = "  -2.90%  "

$ws.Range("D26").Value = "'138.63"
$ws.Range("E26").Value = "  -1.87%  "

$ws.Range("E27").Value = "  -4.02%  "

$ws.Range("E28").Value = "  +1.64%  "

$ws.Range("D29").Value = "'15.58"
$ws.Range("E29").Value = "  +1.02%  "

$ws.Range("E30").Value = "  +0.69%  "

$ws.Range("E31").Value = "  +1.65%  "

$ws.Range("E32").Value = "  +0.56%  "

$ws.Range("D33").Value = "'3.25"
$ws.Range("E33").Value = "  +1.93%  "

$ws.Range("E34").Value = "  +1.28%  "

$ws.Range("E35").Value = "  +0.81%  "

$ws.Range("E36").Value = "  +1.50%  "

$ws.Range("D37").Value = "'2.58"
$ws.Range("E37").Value = "  +1.78%  "

$ws.Range("E38").Value = "  +0.28%  "

$ws.Range("D39").Value = "1.126.03"
$ws.Range("E39").Value = "  -0.19%  "

$ws.Range("E40").Value = "  +0.93%  "

$ws.Range("E41").Value = "  +0.35%  "

$ws.Range("E42").Value = "  -1.62%  "

$ws.Range("D43").Value = "'99.60"
$ws.Range("E43").Value = "  +2.42%  "

$ws.Range("D44").Value = "'0.800"
$ws.Range("E44").Value = "  +0.84%  "

$ws.Range("D45").Value = "0.0₆0108"
$ws.Range("E45").Value = "  -3.11%  "

$ws.Range("D46").Value = "'55.57"
$ws.Range("E46").Value = "  +1.47%  "

$ws.Range("E48").Value = "  -0.44%  "

$ws.Range("D49").Value = "'7.65"
$ws.Range("E49").Value = "  +0.77%  "

$ws.Range("E50").Value = "  -0.52%  "

$ws.Range("E51").Value = "  +0.01%  "
